# design_input_onebyone.xlsx — "Documentation and bug fixing"
#
# 1. Background sheet was missing the "dist_param4" column that exists in
#    the designinput sheet (columns must line up: param1..4, decimals,
#    corr_sheet). Insert a new column before the old "decimals" column (F)
#    and label it "dist_param4". Everything from the old column F onward
#    shifts right by one.
# 2. Column C ("dist_name") was too wide; narrow it down. The newly
#    inserted column gets the same (narrower) width.
# 3. Comment on background!B1 is cleaned up / corrected (typos, authoring
#    note removed, runs merged).
# 4. The "background" sheet becomes the active sheet/tab (was
#    "general_input").

$wb = $excel.ActiveWorkbook

$wsGeneral    = $wb.Worksheets.Item("general_input")
$wsBackground = $wb.Worksheets.Item("background")

# --- 1. Insert the missing "dist_param4" column into background --------
$wsBackground.Columns("F").Insert()
$wsBackground.Range("F1").Value = "dist_param4"

# --- 2. Column widths ----------------------------------------------------
$wsBackground.Columns("C").ColumnWidth = 10.6666666666667
$wsBackground.Columns("F").ColumnWidth = 10.6666666666667

# --- 3. Fix up the B1 cell comment ---------------------------------------
$comment = $wsBackground.Range("B1").Comment
$openQuote = [char]0x201C
$closeQuote = [char]0x201D
$newCommentText = "The set of background parameters is sampled from these distributions. " + "`n" + `
"The number of samples = max(general_input. repeats , designinput.numreals)" + "`n" + `
"The table of sampled values are kept the same for all single sensitivities." + "`n" + `
"This means that for the sensitivities where a parameter is not " + $openQuote + "in focus" + $closeQuote + "  the parameter values will be taken from the background table if it exists there.  If parameter does not exist in background table , parameter defaultvalues are used instead.Distributions and order of distribution parameters are the same as in design input spreadsheet."
$comment.Text($newCommentText)

# --- 4. Make "background" the active sheet/tab, with O16 selected -------
$wsBackground.Activate()
$wsBackground.Range("O16").Select()
